# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the
# f7f39c49-5320-4293-8d9a-99ea6278c1ef row on each sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G is "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 18:44:03"

# zh-cn sheet: column H is "Correspond Handoff Datetime", column K is "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-28 18:43:57"
$wsZhCn.Range("K4").Value = "2016-08-28 18:44:26"

# de-de sheet: column H is "Correspond Handoff Datetime", column K is "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-28 18:44:03"
$wsDeDe.Range("K4").Value = "2016-08-28 18:44:33"
